# Populate Sheet1 with a 4-row x 3-column numeric grid (A1:C4), matching
# the data entered in the target workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1
$ws.Range("B1").Value = 2
$ws.Range("C1").Value = 3

$ws.Range("A2").Value = 4
$ws.Range("B2").Value = 5
$ws.Range("C2").Value = 6

$ws.Range("A3").Value = 7
$ws.Range("B3").Value = 8
$ws.Range("C3").Value = 9

$ws.Range("A4").Value = 10
$ws.Range("B4").Value = 11
$ws.Range("C4").Value = 12

# Leave the selection on the last entered cell, like a user typing down the
# grid and landing on C4.
$ws.Range("C4").Select()

# Target print setup: A4 paper, portrait orientation.
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
